$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.573.72'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').Value = '2.770.01'
$ws.Range('E3').Value = '  -2.57%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''358.62'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').Value = '''108.21'
$ws.Range('E6').Value = '  -4.08%  '
$ws.Range('D7').Value = '''0.552'
$ws.Range('E7').Value = '  -3.30%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = '''0.587'
$ws.Range('E9').Value = '  -3.72%  '
$ws.Range('D10').Value = '''39.51'
$ws.Range('E10').Value = '  -4.12%  '
$ws.Range('E11').Value = '  +4.33%  '
$ws.Range('D12').Value = '''0.0841'
$ws.Range('E12').Value = '  -3.19%  '
$ws.Range('D13').Value = '''19.63'
$ws.Range('E13').Value = '  -2.55%  '
$ws.Range('D14').Value = '''7.59'
$ws.Range('E14').Value = '  -3.26%  '
$ws.Range('D15').Value = '3.210.40'
$ws.Range('E15').Value = '  -2.32%  '
$ws.Range('D16').Value = '2.793.46'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '''0.923'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '51.562.24'
$ws.Range('E18').Value = '  -1.43%  '
$ws.Range('D19').Value = '''7.64'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D20').Value = '''3.08'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('D21').Value = '''13.11'
$ws.Range('E21').Value = '  -2.92%  '
$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('E22').Value = '  -3.73%  '
$ws.Range('D23').Value = '''69.87'
$ws.Range('E23').Value = '  -1.15%  '
$ws.Range('D24').Value = '''267.63'
$ws.Range('E24').Value = '  -2.11%  '
$ws.Range('D25').Value = '''2.76'
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('D26').Value = '''26.24'
$ws.Range('E26').Value = '  -2.93%  '
$ws.Range('D28').Value = '''0.164'
$ws.Range('E28').Value = '  +13.96%  '
$ws.Range('D29').Value = '''10.13'
$ws.Range('E29').Value = '  -2.39%  '
$ws.Range('D30').Value = '''2.25'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').Value = '''35.16'
$ws.Range('E31').Value = '  -1.11%  '
$ws.Range('B32').Value = 'OKB'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D32').Value = '''51.81'
$ws.Range('E32').Value = '  -1.49%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''6.07'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('D34').Value = '''0.0438'
$ws.Range('E34').Value = '  -10.18%  '
$ws.Range('D35').Value = '''0.0838'
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('D36').Value = '''5.14'
$ws.Range('E36').Value = '  -7.83%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = '''18.68'
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('D39').Value = '''3.13'
$ws.Range('E39').Value = '  -5.01%  '
$ws.Range('E40').Value = '  -5.26%  '
$ws.Range('E41').Value = '  -3.31%  '
$ws.Range('D42').Value = '''2.49'
$ws.Range('E42').Value = '  -2.77%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = '''120.38'
$ws.Range('E43').Value = '  -5.04%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '''2.20'
$ws.Range('E44').Value = '  -3.52%  '
$ws.Range('D45').Value = '''21.74'
$ws.Range('E45').Value = '  -5.82%  '
$ws.Range('D46').Value = '2.082.11'
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('D47').Value = '''3.24'
$ws.Range('E47').Value = '  -3.73%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').Value = '''0.927'
$ws.Range('E49').Value = '  -4.92%  '
$ws.Range('D50').Value = '''5.53'
$ws.Range('E50').Value = '  -6.82%  '
$ws.Range('E51').Value = '  +3.94%  '
